$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural change -----------------------------------------------------
# A new "Loco 0 back connection to TP" column is introduced between the
# existing "Loco 0 connected to TP" column (D, renamed to "... front
# connection to TP") and the "TP 0 moveable flag" column. Insert a new
# column at E so the old column E ("TP 0 moveable flag" + its boolean data)
# shifts right to F, carrying its style/type along with it.
$ws.Columns.Item(5).Insert()

# --- Header row --------------------------------------------------------
$ws.Range("D1").Value = "Loco 0 front connection to TP"
$ws.Range("E1").Value = "Loco 0 back connection to TP"
$ws.Range("F1").Value = "TP 0 moveable flag"

# --- Data rows -----------------------------------------------------------
# Columns: A = loco index, B = Loco 0 pos, C = TP 0 pos,
#          D = Loco 0 front connection to TP, E = Loco 0 back connection to TP,
#          F = TP 0 moveable flag (boolean)
$data = @(
    @{ Row = 2;  B = 1;  C = 6; D = -1; E = -1; F = $false }
    @{ Row = 3;  B = 2;  C = 6; D = -1; E = -1; F = $false }
    @{ Row = 4;  B = 3;  C = 6; D = -1; E = -1; F = $false }
    @{ Row = 5;  B = 4;  C = 6; D = -1; E = -1; F = $false }
    @{ Row = 6;  B = 5;  C = 6; D = -1; E = 0;  F = $false }
    @{ Row = 7;  B = 5;  C = 6; D = -1; E = 0;  F = $true }
    @{ Row = 8;  B = 4;  C = 5; D = -1; E = 0;  F = $true }
    @{ Row = 9;  B = 3;  C = 4; D = -1; E = 0;  F = $true }
    @{ Row = 10; B = 4;  C = 5; D = -1; E = 0;  F = $false }
    @{ Row = 11; B = 3;  C = 5; D = -1; E = -1; F = $false }
    @{ Row = 12; B = 2;  C = 5; D = -1; E = -1; F = $false }
    @{ Row = 13; B = 1;  C = 5; D = -1; E = -1; F = $false }
    @{ Row = 14; B = 1;  C = 5; D = -1; E = -1; F = $false }
    @{ Row = 15; B = 0;  C = 5; D = -1; E = -1; F = $false }
    @{ Row = 16; B = 1;  C = 5; D = -1; E = -1; F = $false }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
}
